$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.250.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.22%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.928.25'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.89%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.29'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.96%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7185'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -11.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.18%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3249'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.89%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.43'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.74%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06824'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.77%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8008'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.90%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07938'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.21%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.925.75'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.99%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.39'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.98%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.50'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.87%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '260.34'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.13%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '30.257.57'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.20%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007946'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.08%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.825'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.49%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.180.74'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.88%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9988'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.855'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.31%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.650'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.78%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.98'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.51%  '

# Row 27
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1330'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -11.26%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.93'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.16%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.284'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.362'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.546'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.420'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.35%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.187'

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05068'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.71%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.193'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.89%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7406'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.37%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.730'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.47%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01932'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.805'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.63%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.82'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.565'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.44%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4449'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.36%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.999'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.14%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.25%  '

# Row 45
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8310'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.57%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.71'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.67%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.704'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.56%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.270'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.17%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.14'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.64%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.481'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.23%  '

# Row 51
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4104'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.64%  '
